# Incorporate updated data from upstream processes through 2024
#
# The diff updates the "Open year" = 2024 row (row 26) of Sheet1:
#   - Energy Storage (column C) : 0      -> 7.68
#   - Solar          (column E) : 33.932 -> 74.607
#
# (The workbook's embedded chart caches the same two numbers for plotting;
# updating the worksheet is the authoritative data edit that the chart
# reads from and redraws against on recalculation.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C26").Value = 7.68
$ws.Range("E26").Value = 74.607

# Ask Excel to recalculate / refresh the chart so it picks up the new values.
$excel.CalculateFullRebuild()
$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh()
